$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: number the 17 simulation rows 1..17 (rows 4-20) ---
$bValues = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# --- Recolor the "finished" rows (Jorge row 3, Antonio rows 4-8, Pablo rows 15-20)
#     from the old plain green (RGB 00B050) to the theme "Green, Accent 6" swatch ---
$greenRows = @(3,4,5,6,7,8,15,16,17,18,19,20)
foreach ($row in $greenRows) {
    $r = $ws.Range("C" + $row + ":F" + $row)
    $r.Interior.ThemeColor = 10   # xlThemeColorAccent6 -> OOXML theme index 9 (Accent 6 / green)
    $r.Interior.TintAndShade = 0
}

# --- Move the active selection to I9 ---
$ws.Range("I9").Select()
